$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.637.84"
$ws.Range("E2").Value = "  -4.43%  "
$ws.Range("D3").Value = "3.448.93"
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'599.48"
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("D6").Value = "'147.13"
$ws.Range("E6").Value = "  -7.64%  "
$ws.Range("D7").Value = "3.446.82"
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("E10").Value = "  -5.35%  "
$ws.Range("D11").Value = "'7.38"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'0.422"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("E13").Value = "  -8.05%  "
$ws.Range("D14").Value = "'31.34"
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("D15").Value = "4.017.94"
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("D16").Value = "3.429.82"
$ws.Range("E16").Value = "  -5.21%  "
$ws.Range("D17").Value = "66.638.49"
$ws.Range("E17").Value = "  -4.57%  "
$ws.Range("D18").Value = "'0.117"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").Value = "'9.89"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").Value = "'436.51"
$ws.Range("E22").Value = "  -6.16%  "
$ws.Range("D23").Value = "'0.611"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("D24").Value = "'78.27"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "3.579.26"
$ws.Range("E26").Value = "  -4.89%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("E27").Value = "  -12.67%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = "  -9.54%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  -10.92%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.46"
$ws.Range("E30").Value = "  -6.68%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.59"
$ws.Range("E31").Value = "  -8.43%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.164"
$ws.Range("E33").Value = "  -7.75%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'25.26"
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.07"
$ws.Range("E35").Value = "  -8.22%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.432.05"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  -8.41%  "
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'7.83"
$ws.Range("E39").Value = "  -8.03%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'174.37"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  -9.60%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0880"
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'5.34"
$ws.Range("E44").Value = "  -6.42%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.876"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'29.27"
$ws.Range("E46").Value = "  -9.43%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'46.04"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.22"
$ws.Range("E48").Value = "  -12.02%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.44"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.42"
$ws.Range("E50").Value = "  -12.48%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.980"
$ws.Range("E51").Value = "  -6.37%  "
